{"js": "// Insert a new paragraph right after the \"Include descriptive statistics...\"\n// paragraph (in the Results section) that tells authors how to number and\n// cite tables/figures, per the commit \"figure and table numbering\".\n\nconst anchorText = \"Include descriptive statistics, results from modelling.\";\n\nconst results = context.document.body.search(anchorText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor paragraph not found: \" + anchorText);\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\n\n// Insert a brand-new (empty) paragraph right after the anchor, then fill it\n// with the four text runs used in the source edit.\nconst newParagraph = anchorParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\nnewParagraph.insertText(\n  \"All tables and figures should be numbered, titled, and cited. Use \\u201c\",\n  Word.InsertLocation.end\n);\nnewParagraph.insertText(\"References -> \", Word.InsertLocation.end);\nnewParagraph.insertText(\"Insert caption\\u201d\", Word.InsertLocation.end);\nnewParagraph.insertText(\n  \" to specify your figure and table numbering.\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph right after the \"Include descriptive statistics...\"\n# paragraph (in the Results section) that tells authors how to number and\n# cite tables/figures, per the commit \"figure and table numbering\".\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Include descriptive statistics, results from modelling.\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.MatchCase = $true\n$rng.Find.Execute($anchorText)\n\nif ($rng.Find.Found) {\n    # Collapse the found range to its end, then insert a new paragraph mark\n    # followed by the new sentence.\n    $rng.Collapse(0)\n\n    $openQuote = [char]0x201C\n    $closeQuote = [char]0x201D\n\n    $newText = \"All tables and figures should be numbered, titled, and cited. Use \" + $openQuote + \"References -> Insert caption\" + $closeQuote + \" to specify your figure and table numbering.\"\n\n    $rng.InsertAfter(\"`r\" + $newText)\n}\n"}
